$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.05045533333333333
$ws.Range("H2").Value = 0.151366
$ws.Range("I2").Value = 0.004442474524580737
$ws.Range("J2").Value = 0.004442474524580737
$ws.Range("M2").Value = 1.522526333333333
$ws.Range("N2").Value = 4.567579
$ws.Range("O2").Value = 0.2115373313282365
$ws.Range("P2").Value = 0.2115373313282365
$ws.Range("Q2").Value = 0.07681957365711112
$ws.Range("R2").Value = 0.6913761629140001
$ws.Range("S2").Value = 0.0009397492054234853
$ws.Range("T2").Value = 0.0009397492054234853

# Row 3
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.05045533333333333
$ws.Range("H3").Value = 0.151366
$ws.Range("I3").Value = 0.004442474524580737
$ws.Range("J3").Value = 0.004442474524580737
$ws.Range("O3").Value = 0.4376697219060474
$ws.Range("P3").Value = 0.4376697219060474
$ws.Range("Q3").Value = 0.1589393287148889
$ws.Range("R3").Value = 1.430453958434
$ws.Range("S3").Value = 0.001944336589747951
$ws.Range("T3").Value = 0.001944336589747951

# Row 4
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.05045533333333333
$ws.Range("H4").Value = 0.151366
$ws.Range("I4").Value = 0.004442474524580737
$ws.Range("J4").Value = 0.004442474524580737
$ws.Range("M4").Value = 2.524809666666667
$ws.Range("N4").Value = 7.574429
$ws.Range("O4").Value = 0.3507929467657161
$ws.Range("P4").Value = 0.3507929467657162
$ws.Range("Q4").Value = 0.1273901133348889
$ws.Range("R4").Value = 1.146511020014
$ws.Range("S4").Value = 0.0015583887294093
$ws.Range("T4").Value = 0.001558388729409301

# Row 5
$ws.Range("I5").Value = 0.7425623198471305
$ws.Range("J5").Value = 0.7425623198471305
$ws.Range("M5").Value = 1.522526333333333
$ws.Range("N5").Value = 4.567579
$ws.Range("O5").Value = 0.2115373313282365
$ws.Range("P5").Value = 0.2115373313282365
$ws.Range("Q5").Value = 12.84043847834456
$ws.Range("R5").Value = 115.563946305101
$ws.Range("S5").Value = 0.1570796514853664
$ws.Range("T5").Value = 0.1570796514853664

# Row 6
$ws.Range("I6").Value = 0.7425623198471305
$ws.Range("J6").Value = 0.7425623198471305
$ws.Range("O6").Value = 0.4376697219060474
$ws.Range("P6").Value = 0.4376697219060474
$ws.Range("S6").Value = 0.324997044025403
$ws.Range("T6").Value = 0.324997044025403

# Row 7
$ws.Range("I7").Value = 0.7425623198471305
$ws.Range("J7").Value = 0.7425623198471305
$ws.Range("M7").Value = 2.524809666666667
$ws.Range("N7").Value = 7.574429
$ws.Range("O7").Value = 0.3507929467657161
$ws.Range("P7").Value = 0.3507929467657162
$ws.Range("Q7").Value = 21.29333495558344
$ws.Range("R7").Value = 191.640014600251
$ws.Range("S7").Value = 0.2604856243363611
$ws.Range("T7").Value = 0.2604856243363611

# Row 8
$ws.Range("G8").Value = 2.873389
$ws.Range("H8").Value = 8.620167
$ws.Range("I8").Value = 0.2529952056282888
$ws.Range("J8").Value = 0.2529952056282888
$ws.Range("M8").Value = 1.522526333333333
$ws.Range("N8").Value = 4.567579
$ws.Range("O8").Value = 0.2115373313282365
$ws.Range("P8").Value = 0.2115373313282365
$ws.Range("Q8").Value = 4.374810418410333
$ws.Range("R8").Value = 39.37329376569301
$ws.Range("S8").Value = 0.05351793063744665
$ws.Range("T8").Value = 0.05351793063744665

# Row 9
$ws.Range("G9").Value = 2.873389
$ws.Range("H9").Value = 8.620167
$ws.Range("I9").Value = 0.2529952056282888
$ws.Range("J9").Value = 0.2529952056282888
$ws.Range("O9").Value = 0.4376697219060474
$ws.Range("P9").Value = 0.4376697219060474
$ws.Range("Q9").Value = 9.051461731103668
$ws.Range("R9").Value = 81.46315557993302
$ws.Range("S9").Value = 0.1107283412908964
$ws.Range("T9").Value = 0.1107283412908964

# Row 10
$ws.Range("G10").Value = 2.873389
$ws.Range("H10").Value = 8.620167
$ws.Range("I10").Value = 0.2529952056282888
$ws.Range("J10").Value = 0.2529952056282888
$ws.Range("M10").Value = 2.524809666666667
$ws.Range("N10").Value = 7.574429
$ws.Range("O10").Value = 0.3507929467657161
$ws.Range("P10").Value = 0.3507929467657162
$ws.Range("Q10").Value = 7.254760323293667
$ws.Range("R10").Value = 65.292842909643
$ws.Range("S10").Value = 0.08874893369994571
$ws.Range("T10").Value = 0.08874893369994573
